# Level Renderer parsing Pt 1 - Day 5
# Fill in earned values for the first two CORE FEATURES rows (C4, C5),
# then move the active selection to C8, matching the author's review flow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 3. Run-time file I/O parsing of exported level information (GameLevel.txt)
$ws.Range("C4").Value = 0.05

# 4. Reading binary model data for all referenced models (*.h2b + h2bParser.h)
$ws.Range("C5").Value = 0.01

# Move selection to C8, the next item to grade.
$ws.Range("C8").Select()
